# "Bene Addition and Verification"
# Adds 7 new columns (J:P) of data to Sheet1 describing cheque-book
# transaction-verification lookups (success/status messages + SQL
# queries used to read back the logged transaction), for both the
# header row (row 1) and the sample/case row (row 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 1: new column headers (J1:P1) ----
$ws.Range("J1").Value = "success_message"
$ws.Range("K1").Value = "tran_type_query"
$ws.Range("L1").Value = "tran_date_query"
$ws.Range("M1").Value = "tran_account_no_query"
$ws.Range("N1").Value = "tran_response_query"
$ws.Range("O1").Value = "tran_cheque_no_query"
$ws.Range("P1").Value = "status_message"

# ---- Row 2: new column values (J2:P2) ----
$ws.Range("J2").Value = "Your Cheque Book request has been logged successfully. You can check the status in My Account-->Cheque Book-->Status"

# K2 carries an explicit Text number format (numFmtId 49) in the source
# workbook, so apply that before writing the value.
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "SELECT AC.DESCRIPTION FROM DC_TRANSACTION_ACTIVITY_CONFIG AC WHERE AC.TRANSACTION_TYPE_ID =(SELECT DT.TRANSACTION_TYPE_ID FROM DC_TRANSACTION DT WHERE DT.TRANSACTION_ID = '"

$ws.Range("L2").Value = "SELECT CREATED_ON FROM DC_TRANSACTION DT where DT.TRANSACTION_ID='"
$ws.Range("M2").Value = "SELECT FROM_ACCOUNT FROM DC_TRANSACTION DT where DT.TRANSACTION_ID='"
$ws.Range("N2").Value = "SELECT RESPONSE_MESSAGE FROM DC_TRANSACTION DT where DT.TRANSACTION_ID='"
$ws.Range("O2").Value = "SELECT CHQ_NO_OF_LEAVES FROM DC_TRANSACTION DT where DT.TRANSACTION_ID='"
$ws.Range("P2").Value = "Your Cheque Book Request is in process"

# ---- Column widths to fit the new (much wider) content, matching the
# "best fit" auto-sized widths the author's Excel session produced. ----
$ws.Columns.Item(10).ColumnWidth = 110.02213541666667  # J
$ws.Columns.Item(11).ColumnWidth = 185.02213541666666  # K
$ws.Columns.Item(12).ColumnWidth = 71.87760416666667   # L
$ws.Columns.Item(13).ColumnWidth = 75.30729166666667   # M
$ws.Columns.Item(14).ColumnWidth = 78.87760416666667   # N
$ws.Columns.Item(15).ColumnWidth = 79.02213541666667   # O
$ws.Columns.Item(16).ColumnWidth = 36.022135416666664  # P

# ---- Selection state: the saved view selected the whole sheet
# (A1:XFD1048576) while leaving H12 as the last active cell. ----
$ws.Range("H12").Select() | Out-Null
$ws.Cells.Select() | Out-Null
